$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N32").ClearContents()
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0

$ws.Range("H38").Value = 4699.769
$ws.Range("I38").Value = 3675.875
$ws.Range("J38").Value = 6338
$ws.Range("K38").Value = 11027.625
$ws.Range("L38").Value = 19014
$ws.Range("M38").Value = -10655.625
$ws.Range("N38").Value = -19758

$ws.Range("H74").Value = 3500
$ws.Range("I74").Value = 3500
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3500
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -2564

$ws.Range("H77").Value = 3500
$ws.Range("I77").Value = 3500
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 17500
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -12820

$ws.Range("H88").Value = 972.8
$ws.Range("I88").Value = 1019.8
$ws.Range("J88").Value = 925.8
$ws.Range("K88").Value = 1019.8
$ws.Range("L88").Value = 925.8
$ws.Range("M88").Value = -613.8
$ws.Range("N88").Value = -1737.8

$ws.Range("H91").Value = 972.8
$ws.Range("I91").Value = 1019.8
$ws.Range("J91").Value = 925.8
$ws.Range("K91").Value = 1019.8
$ws.Range("L91").Value = 925.8
$ws.Range("M91").Value = 384.2
$ws.Range("N91").Value = -3733.8

$ws.Range("M98").ClearContents()
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0

$ws.Range("H112").Value = 1799
$ws.Range("I112").Value = 1435
$ws.Range("J112").Value = 1981
$ws.Range("K112").Value = 4305
$ws.Range("L112").Value = 5943
$ws.Range("M112").Value = -3197
$ws.Range("N112").Value = -8159

$ws.Range("M122").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0

$ws.Range("H132").Value = 10043.518
$ws.Range("I132").Value = 9019.639999999999
$ws.Range("J132").Value = 16442.75
$ws.Range("K132").Value = 27058.92
$ws.Range("L132").Value = 49328.25
$ws.Range("M132").Value = -24528.92
$ws.Range("N132").Value = -54388.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N130").ClearContents()
$ws.Range("H130").Value = 14999
$ws.Range("I130").Value = 14999
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 14999
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = -9979

$ws.Range("H132").Value = 2623.1482
$ws.Range("I132").Value = 2469.76
$ws.Range("J132").Value = 4540.5
$ws.Range("K132").Value = 7409.280000000001
$ws.Range("L132").Value = 13621.5
$ws.Range("M132").Value = -4879.280000000001
$ws.Range("N132").Value = -18681.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2955.875
$ws.Range("I20").Value = 2069
$ws.Range("J20").Value = 3488
$ws.Range("K20").Value = 2069
$ws.Range("L20").Value = 3488
$ws.Range("M20").Value = -1822
$ws.Range("N20").Value = -3982

$ws.Range("H22").Value = 1462.5
$ws.Range("I22").Value = 1462.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1462.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1289.5

$ws.Range("H24").Value = 1218.6666
$ws.Range("I24").Value = 1153.5
$ws.Range("J24").Value = 1349
$ws.Range("K24").Value = 1153.5
$ws.Range("L24").Value = 1349
$ws.Range("M24").Value = -918.5
$ws.Range("N24").Value = -1819

$ws.Range("H25").Value = 2128.5
$ws.Range("I25").Value = 2304.6667
$ws.Range("J25").Value = 1600
$ws.Range("K25").Value = 2304.6667
$ws.Range("L25").Value = 1600
$ws.Range("M25").Value = -2069.6667
$ws.Range("N25").Value = -2070

$ws.Range("H29").Value = 682.3333
$ws.Range("I29").Value = 458.8
$ws.Range("J29").Value = 1800
$ws.Range("K29").Value = 458.8
$ws.Range("L29").Value = 1800
$ws.Range("M29").Value = -169.8
$ws.Range("N29").Value = -2378

$ws.Range("M31").ClearContents()
$ws.Range("H31").Value = 2000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 2000
$ws.Range("N31").Value = -2504

$ws.Range("H34").Value = 2577.5
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 2577.5
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 2577.5
$ws.Range("N34").Value = -2805.5

$ws.Range("H36").Value = 109.333336
$ws.Range("I36").Value = 109.333336
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 109.333336
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 424.666664

$ws.Range("H37").Value = 2094.9
$ws.Range("I37").Value = 1389.8
$ws.Range("J37").Value = 2800
$ws.Range("K37").Value = 1389.8
$ws.Range("L37").Value = 2800
$ws.Range("M37").Value = -1252.8
$ws.Range("N37").Value = -3074

$ws.Range("H39").Value = 25000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 25000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 25000
$ws.Range("N39").Value = -25778

$ws.Range("H86").Value = 5923.4736
$ws.Range("I86").Value = 4895.8
$ws.Range("J86").Value = 7065.3335
$ws.Range("K86").Value = 4895.8
$ws.Range("L86").Value = 7065.3335
$ws.Range("M86").Value = -3772.8
$ws.Range("N86").Value = -9311.333500000001

$ws.Range("H89").Value = 5923.4736
$ws.Range("I89").Value = 4895.8
$ws.Range("J89").Value = 7065.3335
$ws.Range("K89").Value = 24479
$ws.Range("L89").Value = 35326.6675
$ws.Range("M89").Value = -18863
$ws.Range("N89").Value = -46558.6675

$ws.Range("H134").Value = 3318
$ws.Range("I134").Value = 3318
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9954
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7419

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3833
$ws.Range("I22").Value = 4249.5
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 4249.5
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -3899.5
$ws.Range("N22").Value = -3700

$ws.Range("H134").Value = 1915.1428
$ws.Range("I134").Value = 1915.1428
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5745.428400000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3210.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 685.5
$ws.Range("I112").Value = 663
$ws.Range("J112").Value = 708
$ws.Range("K112").Value = 1989
$ws.Range("L112").Value = 2124
$ws.Range("M112").Value = -881
$ws.Range("N112").Value = -4340

$ws.Range("H140").Value = 1609.4
$ws.Range("I140").Value = 1232.6666
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 3697.9998
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = 1482.0002
$ws.Range("N140").Value = -25360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 13035125
$ws.Range("I3").Value = 15711000
$ws.Range("J3").Value = 5007500
$ws.Range("K3").Value = 15711000
$ws.Range("L3").Value = 5007500
$ws.Range("M3").Value = -15710884
$ws.Range("N3").Value = -5007732

$ws.Range("H70").Value = 1492.5
$ws.Range("I70").Value = 1492.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 1492.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -1222.5

$ws.Range("H73").Value = 1492.5
$ws.Range("I73").Value = 1492.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 1492.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -556.5

$ws.Range("H102").Value = 1681.05
$ws.Range("I102").Value = 1312.3334
$ws.Range("J102").Value = 4999.5
$ws.Range("K102").Value = 1312.3334
$ws.Range("L102").Value = 4999.5
$ws.Range("M102").Value = 309.6666
$ws.Range("N102").Value = -8243.5

$ws.Range("H122").Value = 2429.6667
$ws.Range("I122").Value = 1360.2
$ws.Range("J122").Value = 7777
$ws.Range("K122").Value = 4080.6
$ws.Range("L122").Value = 23331
$ws.Range("M122").Value = -1630.6
$ws.Range("N122").Value = -28231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0

$ws.Range("H19").Value = 3500
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 3500
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 3500
$ws.Range("N19").Value = -3840

$ws.Range("H22").Value = 1066.6666
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -2590

$ws.Range("H27").Value = 1066.6666
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -2214

$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4864

$ws.Range("N41").ClearContents()
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0

$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N7").ClearContents()
$ws.Range("H7").Value = 1000
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -887

$ws.Range("H122").Value = 3764.3333
$ws.Range("I122").Value = 2176
$ws.Range("J122").Value = 5749.75
$ws.Range("K122").Value = 6528
$ws.Range("L122").Value = 17249.25
$ws.Range("M122").Value = -4078
$ws.Range("N122").Value = -22149.25

$ws.Range("H126").Value = 7324.5
$ws.Range("I126").Value = 5301.3335
$ws.Range("J126").Value = 7998.8887
$ws.Range("K126").Value = 15904.0005
$ws.Range("L126").Value = 23996.6661
$ws.Range("M126").Value = -13434.0005
$ws.Range("N126").Value = -28936.6661
